$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Clear the 'value' data (column D, rows 2-15) while keeping the header in D1
$ws.Range("D2:D15").ClearContents()

# Autofit column B now that content didn't change but reflect Excel's recompute of best-fit width
$ws.Columns("B:B").AutoFit()

# Update the active selection to B9
$ws.Range("B9").Select()
